# Add mod_category column (K) with category values for each module row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column K
$ws.Range("K1").Value = "mod_category"

# Category assigned to each module row (row 2 => index 0, row 47 => index 45)
$categories = @(
    "Computer Science",
    "Data Science and Analytics",
    "Mathematics",
    "Mathematics",
    "Computer Science",
    "Data Science and Analytics",
    "Mathematics",
    "Mathematics",
    "Mathematics",
    "Statistics",
    "Statistics",
    "Data Science and Analytics",
    "Data Science and Analytics",
    "Mathematics",
    "Statistics",
    "Data Science and Analytics",
    "Data Science and Analytics",
    "Statistics",
    "Statistics",
    "Data Science and Analytics",
    "Data Science and Analytics",
    "Data Science and Analytics",
    "Mathematics",
    "Mathematics",
    "Mathematics",
    "Statistics",
    "Statistics",
    "Statistics",
    "Statistics",
    "Statistics",
    "Statistics",
    "Statistics",
    "Statistics",
    "Statistics",
    "Computer Science",
    "Computer Science",
    "Computer Science",
    "Data Science and Analytics",
    "Computer Science",
    "Data Science and Analytics",
    "Computer Science",
    "Computer Science",
    "Computer Science",
    "Computer Science",
    "Data Science and Analytics",
    "Mathematics"
)

for ($i = 0; $i -lt $categories.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $categories[$i]
}

# Row heights auto-recalculated by Excel after the new column text wraps;
# pin them to the values Excel computed on save.
$rowHeights = @{
    2 = 238
    4 = 272
    5 = 255
    9 = 170
    11 = 272
    13 = 238
    14 = 136
    19 = 289
    20 = 204
    21 = 272
    22 = 272
}
foreach ($r in $rowHeights.Keys) {
    $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}

# Restore view state (zoom + final selection) to match the saved workbook
$excel.ActiveWindow.Zoom = 115
$ws.Range("K48").Select()
